{"js": "// Add Erin's name and UID as a new \"Participants\" line, right after the\n// \"Xing Meng (u6483085)\" line \u2014 same tab + spaces layout as the other\n// participant entries.\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\nconst target = paragraphs.items.find((p) => p.text.indexOf(\"Xing Meng (u6483085)\") !== -1);\nif (!target) {\n  throw new Error('Could not find the \"Xing Meng (u6483085)\" participant paragraph');\n}\n\ntarget.insertParagraph(\"\\t         Erin Xiong(u6933612)\", \"After\");\nawait context.sync();\n", "ps1": "# Add Erin's name and UID as a new \"Participants\" line, right after the\n# \"Xing Meng (u6483085)\" line \u2014 same tab + spaces layout as the other\n# participant entries.\n$d = $word.ActiveDocument\n\n$target = $null\nforeach ($p in $d.Paragraphs) {\n    if ($p.Range.Text -like \"*Xing Meng (u6483085)*\") {\n        $target = $p\n        break\n    }\n}\n\nif ($target -eq $null) {\n    throw 'Could not find the \"Xing Meng (u6483085)\" participant paragraph'\n}\n\n$r = $target.Range\n$r.Collapse(0)          # wdCollapseEnd\n$r.InsertParagraphAfter()\n\n$newPara = $target.Next()\n$newPara.Range.Text = \"`t         Erin Xiong(u6933612)\"\n"}
